# Auto-generated script to apply phantom-profits data refresh
# Updates static (non-formula) numeric values across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 134.1
$ws.Range("I6").Value = 132.44444
$ws.Range("K6").Value = 397.33332
$ws.Range("M6").Value = -285.33332
$ws.Range("H93").Value = 408666.66
$ws.Range("I93").Value = 30000
$ws.Range("J93").Value = 598000
$ws.Range("K93").Value = 30000
$ws.Range("L93").Value = 598000
$ws.Range("M93").Value = -27504
$ws.Range("N93").Value = -602992
$ws.Range("H106").Value = 905.1539
$ws.Range("I106").Value = 905.1539
$ws.Range("K106").Value = 905.1539
$ws.Range("M106").Value = -274.1539
$ws.Range("H107").Value = 1353.0625
$ws.Range("I107").Value = 1357.6364
$ws.Range("J107").Value = 1343
$ws.Range("K107").Value = 1357.6364
$ws.Range("L107").Value = 1343
$ws.Range("M107").Value = 562.3635999999999
$ws.Range("N107").Value = -5183
$ws.Range("H111").Value = 3596.5
$ws.Range("J111").Value = 3675
$ws.Range("L111").Value = 11025
$ws.Range("N111").Value = -17159
$ws.Range("H112").Value = 2801.25
$ws.Range("J112").Value = 2701.9
$ws.Range("L112").Value = 8105.700000000001
$ws.Range("N112").Value = -10321.7
$ws.Range("H115").Value = 3240.4285
$ws.Range("I115").Value = 3240.4285
$ws.Range("K115").Value = 9721.2855
$ws.Range("M115").Value = -8154.2855
$ws.Range("H127").Value = 2269
$ws.Range("I127").Value = 2269
$ws.Range("K127").Value = 6807
$ws.Range("M127").Value = -1847
$ws.Range("H132").Value = 4474.515
$ws.Range("I132").Value = 4609.3794
$ws.Range("J132").Value = 3496.75
$ws.Range("K132").Value = 13828.1382
$ws.Range("L132").Value = 10490.25
$ws.Range("M132").Value = -11298.1382
$ws.Range("N132").Value = -15550.25
$ws.Range("H137").Value = 2480.4
$ws.Range("I137").Value = 885.0909
$ws.Range("J137").Value = 3733.8572
$ws.Range("K137").Value = 2655.2727
$ws.Range("L137").Value = 11201.5716
$ws.Range("M137").Value = -105.2727
$ws.Range("N137").Value = -16301.5716
$ws.Range("H138").Value = 1973.8334
$ws.Range("J138").Value = 2473.25
$ws.Range("L138").Value = 7419.75
$ws.Range("N138").Value = -17699.75
$ws.Range("H141").Value = 10078.611
$ws.Range("I141").Value = 9622.066000000001
$ws.Range("J141").Value = 12361.333
$ws.Range("K141").Value = 28866.198
$ws.Range("L141").Value = 37083.999
$ws.Range("M141").Value = -23686.198
$ws.Range("N141").Value = -47443.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 8601.538
$ws.Range("I110").Value = 8485.833000000001
$ws.Range("K110").Value = 8485.833000000001
$ws.Range("M110").Value = -6440.833000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 1982.6
$ws.Range("I75").Value = 1982.6
$ws.Range("K75").Value = 1982.6
$ws.Range("M75").Value = -1046.6
$ws.Range("H78").Value = 1982.6
$ws.Range("I78").Value = 1982.6
$ws.Range("K78").Value = 5947.799999999999
$ws.Range("M78").Value = -1267.799999999999
$ws.Range("H103").Value = 19749.5
$ws.Range("J103").Value = 19749.5
$ws.Range("L103").Value = 19749.5
$ws.Range("N103").Value = -22093.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1175.6666
$ws.Range("I16").Value = 1181.6666
$ws.Range("J16").Value = 1173.2667
$ws.Range("K16").Value = 1181.6666
$ws.Range("L16").Value = 1173.2667
$ws.Range("M16").Value = -894.6666
$ws.Range("N16").Value = -1747.2667
$ws.Range("H19").Value = 25000570
$ws.Range("I19").Value = 28571650
$ws.Range("K19").Value = 28571650
$ws.Range("M19").Value = -28571480
$ws.Range("H24").Value = 25000570
$ws.Range("I24").Value = 28571650
$ws.Range("K24").Value = 28571650
$ws.Range("M24").Value = -28571480
$ws.Range("H58").Value = 2552.625
$ws.Range("I58").Value = 2560.1428
$ws.Range("K58").Value = 2560.1428
$ws.Range("M58").Value = -2357.1428
$ws.Range("H96").Value = 19461.25
$ws.Range("J96").Value = 19461.25
$ws.Range("L96").Value = 19461.25
$ws.Range("N96").Value = -24953.25
$ws.Range("H107").Value = 728.5
$ws.Range("I107").Value = 939.9
$ws.Range("K107").Value = 939.9
$ws.Range("M107").Value = 980.1
$ws.Range("H113").Value = 1175.6666
$ws.Range("I113").Value = 1181.6666
$ws.Range("J113").Value = 1173.2667
$ws.Range("K113").Value = 1181.6666
$ws.Range("L113").Value = 1173.2667
$ws.Range("M113").Value = 988.3334
$ws.Range("N113").Value = -5513.2667
$ws.Range("H134").Value = 2771.4736
$ws.Range("I134").Value = 2772.4375
$ws.Range("K134").Value = 8317.3125
$ws.Range("M134").Value = -5782.3125
$ws.Range("H136").Value = 2552.625
$ws.Range("I136").Value = 2560.1428
$ws.Range("K136").Value = 7680.428400000001
$ws.Range("M136").Value = -5130.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1872551.4
$ws.Range("I4").Value = 1068910
$ws.Range("K4").Value = 3206730
$ws.Range("M4").Value = -3206618
$ws.Range("H8").Value = 5053519.5
$ws.Range("I8").Value = 5053519.5
$ws.Range("K8").Value = 15160558.5
$ws.Range("M8").Value = -15160419.5
$ws.Range("H11").Value = 2063.3333
$ws.Range("I11").Value = 2063.3333
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 6189.999899999999
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -6049.999899999999
$ws.Range("N11").ClearContents()
$ws.Range("H68").Value = 3956.3333
$ws.Range("J68").Value = 5520.4
$ws.Range("L68").Value = 16561.2
$ws.Range("N68").Value = -18183.2
$ws.Range("H71").Value = 3956.3333
$ws.Range("J71").Value = 5520.4
$ws.Range("L71").Value = 49683.6
$ws.Range("N71").Value = -57795.6
$ws.Range("H75").Value = 1275.5
$ws.Range("I75").Value = 700.6667
$ws.Range("K75").Value = 2102.0001
$ws.Range("M75").Value = -1104.0001
$ws.Range("H78").Value = 1275.5
$ws.Range("I78").Value = 700.6667
$ws.Range("K78").Value = 6306.0003
$ws.Range("M78").Value = -1314.0003
$ws.Range("H93").Value = 8095
$ws.Range("I93").Value = 1200
$ws.Range("J93").Value = 14990
$ws.Range("K93").Value = 3600
$ws.Range("L93").Value = 44970
$ws.Range("M93").Value = -1728
$ws.Range("N93").Value = -48714
$ws.Range("H131").Value = 1309.25
$ws.Range("I131").Value = 1309.25
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 3927.75
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 1112.25
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("N37").ClearContents()
$ws.Range("H80").Value = 1124.375
$ws.Range("I80").Value = 999.1667
$ws.Range("K80").Value = 999.1667
$ws.Range("M80").Value = -1.166699999999992
$ws.Range("H83").Value = 1124.375
$ws.Range("I83").Value = 999.1667
$ws.Range("K83").Value = 4995.8335
$ws.Range("M83").Value = -3.833499999999731
$ws.Range("H128").Value = 59980
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 4469
$ws.Range("I132").Value = 3711.25
$ws.Range("K132").Value = 11133.75
$ws.Range("M132").Value = -8603.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 6000
$ws.Range("J14").Value = 6000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = -6344
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("H43").Value = 591666.7
$ws.Range("I43").Value = 25000
$ws.Range("J43").Value = 705000
$ws.Range("K43").Value = 25000
$ws.Range("L43").Value = 705000
$ws.Range("M43").Value = -24807
$ws.Range("N43").Value = -705386
$ws.Range("H80").Value = 70000
$ws.Range("J80").Value = 70000
$ws.Range("L80").Value = 70000
$ws.Range("N80").Value = -72246
$ws.Range("H83").Value = 70000
$ws.Range("J83").Value = 70000
$ws.Range("L83").Value = 210000
$ws.Range("N83").Value = -221232
$ws.Range("H93").Value = 1801.8
$ws.Range("I93").Value = 1702.25
$ws.Range("K93").Value = 1702.25
$ws.Range("M93").Value = -454.25
$ws.Range("H132").Value = 2184.238
$ws.Range("J132").Value = 2660
$ws.Range("L132").Value = 7980
$ws.Range("N132").Value = -13040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1727.3636
$ws.Range("I81").Value = 1727.3636
$ws.Range("K81").Value = 3454.7272
$ws.Range("M81").Value = -2393.7272
$ws.Range("H84").Value = 1727.3636
$ws.Range("I84").Value = 1727.3636
$ws.Range("K84").Value = 17273.636
$ws.Range("M84").Value = -11969.636
$ws.Range("H132").Value = 2015.08
$ws.Range("I132").Value = 2048.7273
$ws.Range("K132").Value = 6146.1819
$ws.Range("M132").Value = -3616.1819
$ws.Range("H136").Value = 3160.0833
$ws.Range("I136").Value = 3177.15
$ws.Range("K136").Value = 9531.450000000001
$ws.Range("M136").Value = -6981.450000000001
